$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Detach the old hyperlink from F2 (it is re-created at G2 further down)
$ws.Range("F2").Hyperlinks.Delete()

# 2) Fully remove cells that have no counterpart at all in the new layout
$toClear = $ws.Range("F2,I2,G3,I3,I4,G5,I5,I6,G7,I7,G8,I8")
foreach ($area in $toClear.Areas) {
    $area.Clear()
}

# 3) Give brand-new cells the same plain style already used across the sheet
#    (copy it in from A1 so no extra style/font entries get created)
$newStyleCells = $ws.Range("J1,E2,H2,J2,F3,H3,J3,H4,J4,F5,H5,J5,H6,J6,F7,H7,J7,H8,J8")
foreach ($area in $newStyleCells.Areas) {
    $ws.Range("A1").Copy($area)
}

# 4) Set cell values
# Row 1
$ws.Range("A1").Value = "TestCaseID"
$ws.Range("B1").Value = "TestStepNumber"
$ws.Range("C1").Value = "TestCaseTitle"
$ws.Range("D1").Value = "Description"
$ws.Range("E1").Value = "Identifier"
$ws.Range("F1").Value = "InputLocator"
$ws.Range("G1").Value = "InputData"
$ws.Range("H1").Value = "Action"
$ws.Range("I1").Value = "Comments"
$ws.Range("J1").Value = "Browser"
# Row 2
$ws.Range("A2").Value = 1.0
$ws.Range("B2").Value = 1.0
$ws.Range("C2").Value = "Verify Whether User is able to login to TCI Dev"
$ws.Range("E2").Value = "xpath"
$ws.Range("G2").Value = "http://tcidev-integration.sandbox.cloud.tibco.com`n"
$ws.Range("H2").Value = "openurl"
$ws.Range("J2").Value = "chrome"
# Row 3
$ws.Range("A3").Value = 1.0
$ws.Range("B3").Value = 2.0
$ws.Range("C3").Value = "Verify Whether User is able to login to TCI Dev"
$ws.Range("E3").Value = "xpath"
$ws.Range("F3").Value = "//button[@id='login']"
$ws.Range("H3").Value = "click"
$ws.Range("J3").Value = "chrome"
# Row 4
$ws.Range("A4").Value = 1.0
$ws.Range("B4").Value = 3.0
$ws.Range("C4").Value = "Verify Whether User is able to login to TCI Dev"
$ws.Range("E4").Value = "xpath"
$ws.Range("F4").Value = "//input[@id='email']"
$ws.Range("G4").Value = "gkchaitu277@dispostable.com"
$ws.Range("H4").Value = "enterText"
$ws.Range("J4").Value = "chrome"
# Row 5
$ws.Range("A5").Value = 1.0
$ws.Range("B5").Value = 4.0
$ws.Range("C5").Value = "Verify Whether User is able to login to TCI Dev"
$ws.Range("F5").Value = "//button[@id='next']"
$ws.Range("H5").Value = "click"
$ws.Range("J5").Value = "chrome"
# Row 6
$ws.Range("A6").Value = 1.0
$ws.Range("B6").Value = 5.0
$ws.Range("C6").Value = "Verify Whether User is able to login to TCI Dev"
$ws.Range("F6").Value = "//input[@id='password']"
$ws.Range("G6").Value = "Tibco2018"
$ws.Range("H6").Value = "enterText"
$ws.Range("J6").Value = "chrome"
# Row 7
$ws.Range("A7").Value = 1.0
$ws.Range("B7").Value = 6.0
$ws.Range("C7").Value = "Verify Whether User is able to login to TCI Dev"
$ws.Range("F7").Value = "//button[@id='taLogin']"
$ws.Range("H7").Value = "click"
$ws.Range("J7").Value = "chrome"
# Row 8
$ws.Range("A8").Value = 1.0
$ws.Range("B8").Value = 7.0
$ws.Range("C8").Value = "Verify Whether User is able to login to TCI Dev"
$ws.Range("H8").Value = "closebrowser"
$ws.Range("J8").Value = "chrome"

# 5) Blank-but-styled cells (style already primed above)
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = ""

# 6) Re-create hyperlink at its new anchor G2
$ws.Hyperlinks.Add($ws.Range("G2"), "http://tcidev-integration.sandbox.cloud.tibco.com")

# 7) Column width changes
$ws.Columns.Item(4).ColumnWidth = 12.7
$ws.Columns.Item(5).ColumnWidth = 18.1667
$ws.Columns.Item(6).ColumnWidth = 18.1667
